# Generate Report for Handback
#
# The "cd7adc74-a919-4fb1-acda-d5106de6dfe9" row (row 5) on both the
# zh-cn and de-de sheets just received a handback. Fill in the
# "Latest Target File", "Latest Handback File" (as a hyperlink),
# "Latest Handback DateTime" and "Error Detail" columns for that row,
# and widen the "Error Detail" column so the long message is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2437f5ba9513ff16962088f3533474852eab704/e2e/cd7adc74-a919-4fb1-acda-d5106de6dfe9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55eba8e256b47765e838a487ae04ff3bc0b8be26/e2e/cd7adc74-a919-4fb1-acda-d5106de6dfe9.md."

$handbackFileDisplay = "cd7adc74-a919-4fb1-acda-d5106de6dfe9.md"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P / column 16) so the new long
# message is legible.
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# Latest Target File -> same xlf file referenced by column G
$wsZh.Range("J5").Value = "cd7adc74-a919-4fb1-acda-d5106de6dfe9.82515ca416e2756668375d70f03868cb2ebf0d5b.zh-cn.xlf"

# Latest Handback DateTime
$wsZh.Range("K5").Value = "2016-10-21 00:40:11"

# Error Detail
$wsZh.Range("P5").Value = $errorDetail

# Latest Handback File, as a hyperlink (mirrors the existing I2/I3/I4 links)
$wsZh.Hyperlinks.Add($wsZh.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b2437f5ba9513ff16962088f3533474852eab704/e2e/cd7adc74-a919-4fb1-acda-d5106de6dfe9.md", "", "", $handbackFileDisplay)

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = 39.17

# Latest Target File -> same xlf file referenced by column G
$wsDe.Range("J5").Value = "cd7adc74-a919-4fb1-acda-d5106de6dfe9.82515ca416e2756668375d70f03868cb2ebf0d5b.de-de.xlf"

# Latest Handback DateTime
$wsDe.Range("K5").Value = "2016-10-21 00:40:29"

# Error Detail
$wsDe.Range("P5").Value = $errorDetail

# Latest Handback File, as a hyperlink (mirrors the existing I2/I3/I4 links)
$wsDe.Hyperlinks.Add($wsDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b2437f5ba9513ff16962088f3533474852eab704/e2e/cd7adc74-a919-4fb1-acda-d5106de6dfe9.md", "", "", $handbackFileDisplay)
